$d = $word.ActiveDocument

$para1 = $d.Paragraphs(1)

# --- 1. Replace the placeholder bookmark text in the paragraph's first run ---
# Paragraph 1 originally holds two runs:
#   "**ID__AFFARS_pgi_5335_topic_6__ID**"  (the run we update)
#   " "                                    (a run that is just a single space, to be removed)
$searchRange = $d.Range($para1.Range.Start, $para1.Range.End)
$searchRange.Find.Execute("**ID__AFFARS_pgi_5335_topic_6__ID**", $true, $false, $false, $false, $false,
                           $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5335_007_90__ID**", 2)

# --- 2. Remove the now-orphaned trailing run that contains only a space ---
$para1 = $d.Paragraphs(1)
$pEnd = $para1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
$trailingSpace.Delete()

# --- 3. Update the paragraph's formatting: indent + paragraph border ---
$para1 = $d.Paragraphs(1)
$pPr = $para1.Range.ParagraphFormat
$pPr.LeftIndent = 11.25
$pPr.Borders.DistanceFromTop = 5
$pPr.Borders.DistanceFromLeft = 5
$pPr.Borders.DistanceFromBottom = 5
$pPr.Borders.DistanceFromRight = 5
